# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 271
$wsExhibition.Range("F5").Value = 152
$wsExhibition.Range("F6").Value = 49
$wsExhibition.Range("F7").Value = 267
$wsExhibition.Range("F9").Value = 1980
$wsExhibition.Range("F10").Value = 350
$wsExhibition.Range("F11").Value = 4667

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 271
$wsAll.Range("F7").Value = 152
$wsAll.Range("F8").Value = 49
$wsAll.Range("F9").Value = 267
$wsAll.Range("F13").Value = 1980
$wsAll.Range("F14").Value = 350
$wsAll.Range("F15").Value = 4667
